$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100 (G100=19906)
$ws.Cells.Item(100, 8).Value = 2170.1667
$ws.Cells.Item(100, 9).Value = 2000
$ws.Cells.Item(100, 10).Value = 2204.2
$ws.Cells.Item(100, 11).Value = 2000
$ws.Cells.Item(100, 12).Value = 2204.2
$ws.Cells.Item(100, 13).Value = -1459
$ws.Cells.Item(100, 14).Value = -3286.2

# Row 106 (G106=19903)
$ws.Cells.Item(106, 8).Value = 168236.67
$ws.Cells.Item(106, 9).Value = 168236.67
$ws.Cells.Item(106, 11).Value = 168236.67
$ws.Cells.Item(106, 13).Value = -167605.67

# Row 129 (G129=36115)
$ws.Cells.Item(129, 8).Value = 831.39703
$ws.Cells.Item(129, 9).Value = 469.58334
$ws.Cells.Item(129, 10).Value = 908.9286
$ws.Cells.Item(129, 11).Value = 1408.75002
$ws.Cells.Item(129, 12).Value = 2726.7858
$ws.Cells.Item(129, 13).Value = 3591.24998
$ws.Cells.Item(129, 14).Value = -12726.7858

$ws = $wb.Worksheets.Item("ARM")
# Row 26 (G26=2241)
$ws.Cells.Item(26, 8).Value = 400
$ws.Cells.Item(26, 9).Value = 400
$ws.Cells.Item(26, 11).Value = 400
$ws.Cells.Item(26, 13).Value = -70

# Row 32 (G32=44147)
$ws.Cells.Item(32, 8).Value = 11618.494
$ws.Cells.Item(32, 9).Value = 11177.746
$ws.Cells.Item(32, 10).Value = 13353.9375
$ws.Cells.Item(32, 11).Value = 11177.746
$ws.Cells.Item(32, 12).Value = 13353.9375
$ws.Cells.Item(32, 13).Value = -10890.746
$ws.Cells.Item(32, 14).Value = -13927.9375

# Row 41 (G41=2501)
$ws.Cells.Item(41, 8).Value = 9800.286
$ws.Cells.Item(41, 9).Value = 2150.5
$ws.Cells.Item(41, 11).Value = 2150.5
$ws.Cells.Item(41, 13).Value = -1736.5

# Row 63 (G63=12528)
$ws.Cells.Item(63, 8).Value = 125002120
$ws.Cells.Item(63, 9).Value = 166668340
$ws.Cells.Item(63, 10).Value = 3498
$ws.Cells.Item(63, 11).Value = 166668340
$ws.Cells.Item(63, 12).Value = 3498
$ws.Cells.Item(63, 13).Value = -166667654
$ws.Cells.Item(63, 14).Value = -4870

# Row 66 (G66=12528)
$ws.Cells.Item(66, 8).Value = 125002120
$ws.Cells.Item(66, 9).Value = 166668340
$ws.Cells.Item(66, 10).Value = 3498
$ws.Cells.Item(66, 11).Value = 833341700
$ws.Cells.Item(66, 12).Value = 17490
$ws.Cells.Item(66, 13).Value = -833338268
$ws.Cells.Item(66, 14).Value = -24354

# Row 102 (G102=19945)
$ws.Cells.Item(102, 8).Value = 1348.2
$ws.Cells.Item(102, 9).Value = 1182.5
$ws.Cells.Item(102, 10).Value = 2011
$ws.Cells.Item(102, 11).Value = 1182.5
$ws.Cells.Item(102, 12).Value = 2011
$ws.Cells.Item(102, 13).Value = 439.5
$ws.Cells.Item(102, 14).Value = -5255

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (G16=27691)
$ws.Cells.Item(16, 8).Value = 1106.3684
$ws.Cells.Item(16, 9).Value = 1002.53845
$ws.Cells.Item(16, 10).Value = 1331.3334
$ws.Cells.Item(16, 11).Value = 1002.53845
$ws.Cells.Item(16, 12).Value = 1331.3334
$ws.Cells.Item(16, 13).Value = -715.53845
$ws.Cells.Item(16, 14).Value = -1905.3334

# Row 31 (G31=44023)
$ws.Cells.Item(31, 8).Value = 7250974
$ws.Cells.Item(31, 9).Value = 7161.2
$ws.Cells.Item(31, 10).Value = 12823137
$ws.Cells.Item(31, 11).Value = 7161.2
$ws.Cells.Item(31, 12).Value = 12823137
$ws.Cells.Item(31, 13).Value = -6866.2
$ws.Cells.Item(31, 14).Value = -12823727

# Row 34 (G34=44023)
$ws.Cells.Item(34, 8).Value = 7250974
$ws.Cells.Item(34, 9).Value = 7161.2
$ws.Cells.Item(34, 10).Value = 12823137
$ws.Cells.Item(34, 11).Value = 7161.2
$ws.Cells.Item(34, 12).Value = 12823137
$ws.Cells.Item(34, 13).Value = -6959.2
$ws.Cells.Item(34, 14).Value = -12823541

# Row 50 (G50=1862)
$ws.Cells.Item(50, 8).Value = 8999.4
$ws.Cells.Item(50, 10).Value = 8999.4
$ws.Cells.Item(50, 12).Value = 8999.4
$ws.Cells.Item(50, 14).Value = -10249.4

# Row 68 (G68=10611)
$ws.Cells.Item(68, 8).Value = 17052.857
$ws.Cells.Item(68, 9).Value = 2980
$ws.Cells.Item(68, 10).Value = 19398.334
$ws.Cells.Item(68, 11).Value = 2980
$ws.Cells.Item(68, 12).Value = 19398.334
$ws.Cells.Item(68, 13).Value = -2231
$ws.Cells.Item(68, 14).Value = -20896.334

# Row 71 (G71=10611)
$ws.Cells.Item(71, 8).Value = 17052.857
$ws.Cells.Item(71, 9).Value = 2980
$ws.Cells.Item(71, 10).Value = 19398.334
$ws.Cells.Item(71, 11).Value = 8940
$ws.Cells.Item(71, 12).Value = 58195.00199999999
$ws.Cells.Item(71, 13).Value = -5196
$ws.Cells.Item(71, 14).Value = -65683.00199999999

# Row 80 (G80=12015)
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

# Row 83 (G83=12015)
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

# Row 105 (G105=19928)
$ws.Cells.Item(105, 8).Value = 1820
$ws.Cells.Item(105, 9).Value = 2155
$ws.Cells.Item(105, 10).Value = 1150
$ws.Cells.Item(105, 11).Value = 2155
$ws.Cells.Item(105, 12).Value = 1150
$ws.Cells.Item(105, 13).Value = -408
$ws.Cells.Item(105, 14).Value = -4644

# Row 113 (G113=27691)
$ws.Cells.Item(113, 8).Value = 1106.3684
$ws.Cells.Item(113, 9).Value = 1002.53845
$ws.Cells.Item(113, 10).Value = 1331.3334
$ws.Cells.Item(113, 11).Value = 1002.53845
$ws.Cells.Item(113, 12).Value = 1331.3334
$ws.Cells.Item(113, 13).Value = 1167.46155
$ws.Cells.Item(113, 14).Value = -5671.3334

# Row 134 (G134=44020)
$ws.Cells.Item(134, 8).Value = 628524.6
$ws.Cells.Item(134, 9).Value = 1955.4814
$ws.Cells.Item(134, 10).Value = 2166467
$ws.Cells.Item(134, 11).Value = 5866.4442
$ws.Cells.Item(134, 12).Value = 6499401
$ws.Cells.Item(134, 13).Value = -3331.4442
$ws.Cells.Item(134, 14).Value = -6504471

$ws = $wb.Worksheets.Item("CUL")
# Row 109 (G109=27854)
$ws.Cells.Item(109, 8).Value = 4060
$ws.Cells.Item(109, 10).Value = 4478.5713
$ws.Cells.Item(109, 12).Value = 13435.7139
$ws.Cells.Item(109, 14).Value = -15515.7139

# Row 112 (G112=27855)
$ws.Cells.Item(112, 8).Value = 4740.919
$ws.Cells.Item(112, 9).Value = 3963.5
$ws.Cells.Item(112, 10).Value = 4835.1514
$ws.Cells.Item(112, 11).Value = 11890.5
$ws.Cells.Item(112, 12).Value = 14505.4542
$ws.Cells.Item(112, 13).Value = -10782.5
$ws.Cells.Item(112, 14).Value = -16721.4542

# Row 118 (G118=27872)
$ws.Cells.Item(118, 8).Value = 1643.5416
$ws.Cells.Item(118, 9).Value = 2457.4
$ws.Cells.Item(118, 10).Value = 1429.3684
$ws.Cells.Item(118, 11).Value = 7372.200000000001
$ws.Cells.Item(118, 12).Value = 4288.1052
$ws.Cells.Item(118, 13).Value = -6129.200000000001
$ws.Cells.Item(118, 14).Value = -6774.1052

# Row 121 (G121=27878)
$ws.Cells.Item(121, 8).Value = 613.8570999999999
$ws.Cells.Item(121, 10).Value = 910.1
$ws.Cells.Item(121, 12).Value = 2730.3
$ws.Cells.Item(121, 14).Value = -5350.3

# Row 131 (G131=36060)
$ws.Cells.Item(131, 8).Value = 2855.2693
$ws.Cells.Item(131, 9).Value = 2201.682
$ws.Cells.Item(131, 10).Value = 6450
$ws.Cells.Item(131, 11).Value = 6605.045999999999
$ws.Cells.Item(131, 12).Value = 19350
$ws.Cells.Item(131, 13).Value = -1565.045999999999
$ws.Cells.Item(131, 14).Value = -29430

$ws = $wb.Worksheets.Item("GSM")
# Row 42 (G42=27213)
$ws.Cells.Item(42, 8).Value = 50000
$ws.Cells.Item(42, 10).Value = 50000
$ws.Cells.Item(42, 12).Value = 50000
$ws.Cells.Item(42, 14).Value = -50970

# Row 57 (G57=2876)
$ws.Cells.Item(57, 8).Value = 12414.143
$ws.Cells.Item(57, 10).Value = 23666.666
$ws.Cells.Item(57, 12).Value = 23666.666
$ws.Cells.Item(57, 14).Value = -25306.666

# Row 112 (G112=25859)
$ws.Cells.Item(112, 8).Value = 45000
$ws.Cells.Item(112, 10).Value = 45000
$ws.Cells.Item(112, 12).Value = 45000
$ws.Cells.Item(112, 14).Value = -47216

# Row 113 (G113=27710)
$ws.Cells.Item(113, 8).Value = 40766.2
$ws.Cells.Item(113, 9).Value = 67171.734
$ws.Cells.Item(113, 10).Value = 1157.9
$ws.Cells.Item(113, 11).Value = 67171.734
$ws.Cells.Item(113, 12).Value = 1157.9
$ws.Cells.Item(113, 13).Value = -65001.734
$ws.Cells.Item(113, 14).Value = -5497.9

# Row 115 (G115=27213)
$ws.Cells.Item(115, 8).Value = 50000
$ws.Cells.Item(115, 10).Value = 50000
$ws.Cells.Item(115, 12).Value = 50000
$ws.Cells.Item(115, 14).Value = -52350

# Row 118 (G118=26172)
$ws.Cells.Item(118, 8).Value = 14111.111
$ws.Cells.Item(118, 10).Value = 14111.111
$ws.Cells.Item(118, 12).Value = 14111.111
$ws.Cells.Item(118, 14).Value = -17425.111

# Row 122 (G122=36182)
$ws.Cells.Item(122, 8).Value = 4168442.5
$ws.Cells.Item(122, 9).Value = 5129483.5
$ws.Cells.Item(122, 10).Value = 3932.6667
$ws.Cells.Item(122, 11).Value = 15388450.5
$ws.Cells.Item(122, 12).Value = 11798.0001
$ws.Cells.Item(122, 13).Value = -15386000.5
$ws.Cells.Item(122, 14).Value = -16698.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 100 (G100=19995)
$ws.Cells.Item(100, 8).Value = 2445.4707
$ws.Cells.Item(100, 9).Value = 2126.6365
$ws.Cells.Item(100, 11).Value = 2126.6365
$ws.Cells.Item(100, 13).Value = -1585.6365

# Row 122 (G122=36247)
$ws.Cells.Item(122, 8).Value = 5370.927
$ws.Cells.Item(122, 9).Value = 5340.0835
$ws.Cells.Item(122, 10).Value = 5414.4707
$ws.Cells.Item(122, 11).Value = 16020.2505
$ws.Cells.Item(122, 12).Value = 16243.4121
$ws.Cells.Item(122, 13).Value = -13570.2505
$ws.Cells.Item(122, 14).Value = -21143.4121

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (G81=12596)
$ws.Cells.Item(81, 8).Value = 966.6667
$ws.Cells.Item(81, 10).Value = 900
$ws.Cells.Item(81, 12).Value = 1800
$ws.Cells.Item(81, 14).Value = -3922

# Row 84 (G84=12596)
$ws.Cells.Item(84, 8).Value = 966.6667
$ws.Cells.Item(84, 10).Value = 900
$ws.Cells.Item(84, 12).Value = 9000
$ws.Cells.Item(84, 14).Value = -19608
